# Seat Assignments sheet: move the 3 seat-assignment rows for block 1
# (IDs eef4b984..., 385ae883..., 6bde6343...) that currently sit at rows
# 45:47 down to the bottom of the table (new rows 86:88), shifting every
# row in between up by three. Net effect matches the upstream commit that
# re-ordered these rows (no column/header changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seat Assignments")

# 1) Copy the 3 rows that need to move (full rows, so every populated
#    column A:H and any per-cell formatting travels with them) to a
#    staging area just past the current last row (row 88), i.e. rows
#    89:91.
$ws.Rows("45:47").Copy($ws.Range("A89"))

# 2) Delete the original rows 45:47. This shifts every row below them
#    (including the staged copy we just placed at 89:91) up by three
#    rows, so the staged copy lands exactly at the new bottom rows
#    86:88 - precisely where they need to end up.
$ws.Rows("45:47").Delete()
